$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for "Haba" was recorded. It belongs right
# before the existing row 15, so insert a new row there; this pushes the
# former rows 15-39 down to 16-40 (the worksheet grows from 39 to 40 rows).
$ws.Rows(15).Insert()

# Fill in the newly inserted row 15 with the new data point. The columns
# that are constant across this whole sheet (market/region/product info)
# are repeated, matching the pattern of every other row.
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44495
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 100112026
$ws.Cells.Item(15, 7).Value = "Haba"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 20
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 9000
$ws.Cells.Item(15, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 360
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
